# Week 6 Soccer Analysis
# Slide 4 ("Research Question(s)") — reword the last bullet in the body
# placeholder from:
#   "Is it more profitable to bet on the favourite, draw, underdog?"
# to:
#   "Is it more profitable to bet on the home team, draw, away team?"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# The body placeholder has 5 paragraphs; the bullet we need is the 4th one.
$para = $tr.Paragraphs(4, 1)

# Step 1: rewrite the tail of the sentence ("underdog?" -> "away team?").
# Replacing just this sub-range splits the sentence into two runs: the
# untouched head and this freshly-edited tail.
$tail = $para.Characters(54, 9)
$tail.Text = "away team?"

# Step 2: rewrite the head run's *entire* span ("favourite" -> "home team")
# so it stays a single run instead of fragmenting further.
$head = $para.Characters(1, 53)
$head.Text = "Is it more profitable to bet on the home team, draw, "
